# Apply updates to the WR_89775268_WeekEnding_062925.xlsx report
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info updates ---
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:00 AM"
$ws.Range("C8").Value = 3738.51
$ws.Range("C9").Value = 11
$ws.Range("G10").Value = ""

# --- Thursday (06/26/2025) section pricing updates ---
$ws.Range("H16").Value = 198.88
$ws.Range("H17").Value = 478.55
$ws.Range("H18").Value = 677.4300000000001

# --- Saturday (06/28/2025) section pricing updates ---
$ws.Range("H23").Value = 1616.94
$ws.Range("H24").Value = 94.17
$ws.Range("H25").Value = 94.17
$ws.Range("H26").Value = 478.55
$ws.Range("H27").Value = 55.18

# Remove the "Point 07 / TIE-4-ALH-F" line item row (row 28), shifting all
# following rows up by one. This collapses the 6 "Point 09" line items down
# to the 4 that remain after dedup, and moves the TOTAL row from 33 to 32.
$ws.Rows("28").Delete()

# --- Update the now-shifted line items (rows 28-31) ---
$ws.Range("H28").Value = 94.17

$ws.Range("B29").Value = "PIN-15-PTP-C"
$ws.Range("D29").Value = "Pin,15kV,Pole top,Corrosive"
$ws.Range("H29").Value = 94.17

$ws.Range("B30").Value = "POL-40-2"
$ws.Range("D30").Value = "Pole,40ft,Class 2"
$ws.Range("H30").Value = 478.55

$ws.Range("B31").Value = "SAA-3-CV-C"
$ws.Range("D31").Value = "SAA,3 inch,Clevis,Corr"
$ws.Range("H31").Value = 55.18

# --- TOTAL row (now row 32 after the deletion) ---
$ws.Range("H32").Value = 3061.08
